$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    0.2432490129025711,
    0.2411897279234069,
    0.2279374252521792,
    0.2110671728913582,
    0.2232835019967376,
    0.2338366062972175,
    0.2135195243514688,
    0.200607347540649,
    0.2823875815070244,
    0.1998897476382224,
    0.2053661800737224,
    0.2025752555140319,
    0.1739384242826933,
    0.2015725473051372,
    0.2320889448780388,
    0.2167895105100233,
    0.2144442310725814,
    0.2121829174672534,
    0.1731373828639528,
    0.1940917271092404,
    0.2079469975294282,
    0.179211919645771,
    0.1674113558621942,
    0.1897901160543019,
    0.205923835628804,
    0.2205501952095101,
    0.2019318048492107,
    0.1958749450597683,
    0.1981261388742524,
    0.2391133932550645,
    0.1950598397699665,
    0.1487496675782086,
    0.1892665454765523,
    0.1765637362566932,
    0.2337421648130379,
    0.196159914135432,
    0.2096921815164276,
    0.1988977444286743,
    0.2290684984032781,
    0.2153576239808776,
    0.2127302161003796,
    0.2123330059331108,
    0.2111907239136736,
    0.2413968397084669,
    0.1838156528288608,
    0.2464650140770868,
    0.2234666136291428,
    0.1959046171534852,
    0.2027371075584232,
    0.2207704047530455,
    0.2514823645950497,
    0.2046301203171215,
    0.1956943611878062,
    0.2466881491246491,
    0.2123710190639136,
    0.1820696510099117,
    0.210302979833812,
    0.1794839063464464,
    0.2233297240657419,
    0.2393133728178417,
    0.1629844592048781,
    0.1969827144392482,
    0.1965733061178754,
    0.183836812759245,
    0.1626435010973513,
    0.1995352458803386,
    0.2485047676408502,
    0.2187800333905454,
    0.2045486734926522,
    0.1906019554497276,
    0.2207784478459045,
    0.2263140766664,
    0.209776060412968,
    0.2082701329743925,
    0.2040152875554088,
    0.2138425143858722,
    0.1728289216916878
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = [double]$newValues[$i]
}

$wb.Save()
